# Updates the cryptos price/volume snapshot data (columns D and E)
# to reflect the latest values from the scheduled GitHub Actions refresh.
# Values are written as literal text (Value2) so that formats such as
# "1.00", "0.0258" and "x.xxx.xx"-style thousand separators are preserved
# exactly, matching how the source data feed renders them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "60.925.69"
$ws.Range("E2").Value2 = "  +0.23%  "
$ws.Range("D3").Value2 = "3.372.48"
$ws.Range("E3").Value2 = "  -0.41%  "
$ws.Range("E4").Value2 = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "572.44"
$ws.Range("E5").Value2 = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "138.53"
$ws.Range("E6").Value2 = "  -1.84%  "
$ws.Range("E7").Value2 = "  +0.00%  "
$ws.Range("E8").Value2 = "  -0.45%  "
$ws.Range("E9").Value2 = "  +2.72%  "
$ws.Range("E10").Value2 = "  -1.80%  "
$ws.Range("E11").Value2 = "  -3.30%  "
$ws.Range("D12").Value2 = "3.945.40"
$ws.Range("E12").Value2 = "  -0.48%  "
$ws.Range("E13").Value2 = "  +0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "28.05"
$ws.Range("E14").Value2 = "  -1.50%  "
$ws.Range("D15").Value2 = "3.384.13"
$ws.Range("E15").Value2 = "  -0.16%  "
$ws.Range("E16").Value2 = "  -1.58%  "
$ws.Range("D17").Value2 = "61.018.15"
$ws.Range("E17").Value2 = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "6.09"
$ws.Range("E18").Value2 = "  -1.79%  "
$ws.Range("E19").Value2 = "  -3.28%  "
$ws.Range("E20").Value2 = "  -0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "384.05"
$ws.Range("E21").Value2 = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "74.96"
$ws.Range("E22").Value2 = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.551"
$ws.Range("E23").Value2 = "  -1.43%  "
$ws.Range("E24").Value2 = "  -0.13%  "
$ws.Range("E25").Value2 = "  -5.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "1.00"
$ws.Range("E27").Value2 = "  +0.08%  "
$ws.Range("E28").Value2 = "  -3.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "7.92"
$ws.Range("E29").Value2 = "  -0.75%  "
$ws.Range("E30").Value2 = "  -1.28%  "
$ws.Range("E32").Value2 = "  -6.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "23.02"
$ws.Range("E33").Value2 = "  -2.38%  "
$ws.Range("E34").Value2 = "  -1.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "167.15"
$ws.Range("E35").Value2 = "  +0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "4.91"
$ws.Range("E36").Value2 = "  -1.02%  "
$ws.Range("D37").Value2 = "3.408.05"
$ws.Range("E37").Value2 = "  -0.25%  "
$ws.Range("E38").Value2 = "  -2.72%  "
$ws.Range("E39").Value2 = "  -2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "25.63"
$ws.Range("E40").Value2 = "  -9.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.772"
$ws.Range("E41").Value2 = "  -0.81%  "
$ws.Range("E42").Value2 = "  -1.44%  "
$ws.Range("E43").Value2 = "  -1.78%  "
$ws.Range("E44").Value2 = "  -1.44%  "
$ws.Range("D45").Value2 = "2.447.50"
$ws.Range("E45").Value2 = "  -1.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "6.64"
$ws.Range("E46").Value2 = "  -2.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.00"
$ws.Range("E47").Value2 = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "22.06"
$ws.Range("E48").Value2 = "  -6.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0258"
$ws.Range("E49").Value2 = "  -4.71%  "
$ws.Range("E50").Value2 = "  -4.01%  "
$ws.Range("E51").Value2 = "  -2.96%  "
